$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-text header updates (rich-text runs re-flattened) ---
# A8: "Volume 30   Number  12" -> "Volume 30   Number  13"
$volRng = $ws.Range("A8")
$volRng.Characters(21,2).Text = "13"

# C9: "Report Covering the Week  3/20/2023  Through  3/26/2023"
#  -> "Report Covering the Week  3/27/2023  Through  4/2/2023"
$weekRng = $ws.Range("C9")
$weekRng.Characters(27,9).Text = "3/27/2023"
$weekRng.Characters(47,9).Text = "4/2/2023"

# --- Column E width (bestFit) ---
$ws.Columns("E").ColumnWidth = 8

# --- Updated crime-statistics figures ---
$ws.Range("C15").Value = 1
$ws.Range("I15").Value = 14
$ws.Range("K15").Value = 250
$ws.Range("L15").Value = 600
$ws.Range("M15").Value = 250
$ws.Range("N15").Value = 600
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 150
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 6.25
$ws.Range("I16").Value = 63
$ws.Range("J16").Value = 49
$ws.Range("K16").Value = 28.571428571428
$ws.Range("L16").Value = 162.5
$ws.Range("M16").Value = 28.571428571428
$ws.Range("N16").Value = -79.742765273311
$ws.Range("C17").Value = 5
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 57.142857142857
$ws.Range("I17").Value = 52
$ws.Range("J17").Value = 47
$ws.Range("K17").Value = 10.638297872340
$ws.Range("L17").Value = 13.043478260869
$ws.Range("M17").Value = 62.5
$ws.Range("N17").Value = -17.460317460317
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 70
$ws.Range("I18").Value = 54
$ws.Range("J18").Value = 42
$ws.Range("K18").Value = 28.571428571428
$ws.Range("L18").Value = 45.945945945945
$ws.Range("M18").Value = -23.943661971831
$ws.Range("N18").Value = -84.788732394366
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 63
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = 34.042553191489
$ws.Range("I19").Value = 186
$ws.Range("J19").Value = 171
$ws.Range("K19").Value = 8.771929824561
$ws.Range("L19").Value = 91.752577319587
$ws.Range("M19").Value = 73.831775700934
$ws.Range("N19").Value = 0.540540540540
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 17
$ws.Range("H20").Value = 13.333333333333
$ws.Range("I20").Value = 49
$ws.Range("J20").Value = 41
$ws.Range("K20").Value = 19.512195121951
$ws.Range("L20").Value = 16.666666666666
$ws.Range("M20").Value = -14.035087719298
$ws.Range("N20").Value = -89.876033057851
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = 106.25
$ws.Range("F21").Value = 141
$ws.Range("G21").Value = 103
$ws.Range("H21").Value = 36.893203883495
$ws.Range("I21").Value = 418
$ws.Range("J21").Value = 354
$ws.Range("K21").Value = 18.079096045197
$ws.Range("L21").Value = 68.548387096774
$ws.Range("M21").Value = 29.813664596273
$ws.Range("N21").Value = -70.249110320284
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 400
$ws.Range("F22").Value = 13
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = 85.714285714285
$ws.Range("I22").Value = 25
$ws.Range("J22").Value = 23
$ws.Range("K22").Value = 8.695652173913
$ws.Range("L22").Value = 212.5
$ws.Range("M22").Value = 127.272727272727
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = 38.095238095238
$ws.Range("F24").Value = 123
$ws.Range("G24").Value = 95
$ws.Range("H24").Value = 29.473684210526
$ws.Range("I24").Value = 412
$ws.Range("J24").Value = 290
$ws.Range("K24").Value = 42.068965517241
$ws.Range("L24").Value = 57.854406130268
$ws.Range("M24").Value = 113.471502590674
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 10
$ws.Range("F25").Value = 40
$ws.Range("G25").Value = 43
$ws.Range("H25").Value = -6.976744186046
$ws.Range("I25").Value = 136
$ws.Range("J25").Value = 145
$ws.Range("K25").Value = -6.206896551724
$ws.Range("L25").Value = 51.111111111111
$ws.Range("M25").Value = 3.816793893129
$ws.Range("C26").Value = 1
$ws.Range("I26").Value = 16
$ws.Range("K26").Value = 300
$ws.Range("L26").Value = 220
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 11
$ws.Range("H27").Value = -63.636363636363
$ws.Range("I27").Value = 21
$ws.Range("J27").Value = 23
$ws.Range("K27").Value = -8.695652173913
$ws.Range("L27").Value = 110
$ws.Range("J30").Value = 3
$ws.Range("K30").Value = -33.333333333333

# D30 / E30 (Hate Crimes, Week-to-Date) flip from "N/A" placeholders to real numbers
$ws.Range("D30").Value = 1
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("E30").Value = -100
$ws.Range("E30").NumberFormat = '#,##0.0;"-"#,##0.0'
